# Update the "Förändrad" (Changed) date column (C) for rows 2-72
# from 45175 (2023-09-06) to 45177 (2023-09-08), leaving all other
# cell contents, formatting and styles untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 72; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
